$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Greece Super League 1")

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")

    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()

    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows 36 37
Swap-Rows 124 125
